# Adds a new "Docentes responsáveis:" block (label + two names, each
# duplicated in the B/C "current" / "modified" columns) right after the
# "Objectives:" row, pushing "Programa resumido:" and everything below it
# down by three rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the 3 new rows (old row 12 "Programa resumido:" -> new row 15).
$ws.Rows("12:14").Insert()

# New label (column A, bold style) and the two docentes (columns B & C).
$ws.Range("A12").Value = "Docentes responsáveis:"
$ws.Range("B13").Value = "3577649 - Carlos Angelo Nunes"
$ws.Range("C13").Value = "3577649 - Carlos Angelo Nunes"
$ws.Range("B14").Value = "1922320 - Sebastiao Ribeiro"
$ws.Range("C14").Value = "1922320 - Sebastiao Ribeiro"

# Row-insert only copies column-level formatting; pull in the exact cell
# styles used by the rest of the sheet (A = bold label, B = normal wrap,
# C = red "modified" wrap) from the next data row (now row 15, formerly 12).
$ws.Range("A15").Copy()
$ws.Range("A12").PasteSpecial(-4122)

$ws.Range("B15").Copy()
$ws.Range("B13:B14").PasteSpecial(-4122)

$ws.Range("C15").Copy()
$ws.Range("C13:C14").PasteSpecial(-4122)

# The insert also stamped column-A's bold style onto the now-empty A13/A14
# cells; clear them so those rows only carry the B/C docente cells, matching
# the target layout.
$ws.Range("A13:A14").Clear()

Write-Host "Inserted Docentes responsaveis block (rows 12-14)"
